{"js": "// Insert five new paragraphs at the very start of the document body,\n// before the existing (empty) paragraph, in the order they should read.\nconst body = context.document.body;\n\nconst p1 = body.insertParagraph(\"Queue \", \"Start\");\nconst p2 = p1.insertParagraph(\"Pair-Programming\", \"After\");\nconst p3 = p2.insertParagraph(\"Advisor-Meetings\", \"After\");\nconst p4 = p3.insertParagraph(\"Message Instructor\", \"After\");\nconst p5 = p4.insertParagraph(\"Or Mentor\", \"After\");\n\nawait context.sync();\n", "ps1": "# Insert five new paragraphs at the very start of the document, before\n# the existing (empty) paragraph, in the order they should read.\n$d = $word.ActiveDocument\n\n$r = $d.Paragraphs.First.Range\n$r.InsertBefore(\"Queue `rPair-Programming`rAdvisor-Meetings`rMessage Instructor`rOr Mentor`r\")\n"}
